$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = 4
$ws.Range("B27").Value = "Admin crud"

$ws.Range("J23").Select()
